$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = '@'
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '67.340.70'
Set-TextValue $ws.Range('E2') '  -0.31%  '
Set-TextValue $ws.Range('D3') '2.624.52'
Set-TextValue $ws.Range('E3') '  -2.05%  '
Set-TextValue $ws.Range('D5') '595.34'
Set-TextValue $ws.Range('E5') '  -0.78%  '
Set-TextValue $ws.Range('D6') '166.89'
Set-TextValue $ws.Range('E6') '  +0.69%  '
Set-TextValue $ws.Range('E7') '  +0.02%  '
Set-TextValue $ws.Range('E8') '  -2.36%  '
Set-TextValue $ws.Range('D9') '2.624.39'
Set-TextValue $ws.Range('D10') '0.139'
Set-TextValue $ws.Range('E10') '  -2.33%  '
Set-TextValue $ws.Range('E11') '  +1.15%  '
Set-TextValue $ws.Range('D12') '0.364'
Set-TextValue $ws.Range('E12') '  +1.49%  '
Set-TextValue $ws.Range('E13') '  +0.33%  '
Set-TextValue $ws.Range('D14') '27.65'
Set-TextValue $ws.Range('E14') '  -0.62%  '
Set-TextValue $ws.Range('D15') '3.101.17'
Set-TextValue $ws.Range('E15') '  -2.08%  '
Set-TextValue $ws.Range('E16') '  -1.15%  '
Set-TextValue $ws.Range('D17') '67.246.22'
Set-TextValue $ws.Range('E17') '  -0.32%  '
Set-TextValue $ws.Range('D18') '2.621.25'
Set-TextValue $ws.Range('E18') '  -1.80%  '
Set-TextValue $ws.Range('D19') '12.07'
Set-TextValue $ws.Range('E19') '  +2.65%  '
Set-TextValue $ws.Range('E20') '  +4.04%  '
Set-TextValue $ws.Range('D21') '357.91'
Set-TextValue $ws.Range('E22') '  -1.29%  '
Set-TextValue $ws.Range('D23') '4.67'
Set-TextValue $ws.Range('E23') '  -3.23%  '
Set-TextValue $ws.Range('E24') '  -0.02%  '
Set-TextValue $ws.Range('D25') '1.94'
Set-TextValue $ws.Range('E25') '  -4.96%  '
Set-TextValue $ws.Range('E26') '  +1.00%  '
Set-TextValue $ws.Range('D27') '69.77'
Set-TextValue $ws.Range('E27') '  -1.89%  '
Set-TextValue $ws.Range('D28') '2.759.81'
Set-TextValue $ws.Range('E28') '  -1.65%  '
Set-TextValue $ws.Range('D29') '1.00'
Set-TextValue $ws.Range('E29') '  +0.22%  '
Set-TextValue $ws.Range('E30') '  -2.27%  '
Set-TextValue $ws.Range('D31') '546.34'
Set-TextValue $ws.Range('E31') '  -2.02%  '
Set-TextValue $ws.Range('D32') '7.92'
Set-TextValue $ws.Range('E32') '  -1.14%  '
Set-TextValue $ws.Range('E33') '  -2.96%  '
Set-TextValue $ws.Range('E34') '  -1.77%  '
Set-TextValue $ws.Range('E35') '  +4.43%  '
Set-TextValue $ws.Range('E36') '  +0.04%  '
Set-TextValue $ws.Range('D37') '1.51'
Set-TextValue $ws.Range('E37') '  -3.75%  '
Set-TextValue $ws.Range('D38') '157.01'
Set-TextValue $ws.Range('E38') '  +0.72%  '
Set-TextValue $ws.Range('D39') '19.00'
Set-TextValue $ws.Range('E39') '  -2.90%  '
Set-TextValue $ws.Range('E40') '  -2.14%  '
Set-TextValue $ws.Range('B41') 'RenderToken'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range('D41') '5.21'
Set-TextValue $ws.Range('E41') '  -1.87%  '
Set-TextValue $ws.Range('E42') '  -1.14%  '
Set-TextValue $ws.Range('B43') 'WhiteBITCoin'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D43') '18.18'
Set-TextValue $ws.Range('E43') '  +1.37%  '
Set-TextValue $ws.Range('E44') '  +0.06%  '
Set-TextValue $ws.Range('D45') '2.42'
Set-TextValue $ws.Range('E45') '  -4.45%  '
Set-TextValue $ws.Range('D47') '152.09'
Set-TextValue $ws.Range('E47') '  -0.86%  '
Set-TextValue $ws.Range('D48') '0.580'
Set-TextValue $ws.Range('E48') '  -2.09%  '
Set-TextValue $ws.Range('D49') '3.77'
Set-TextValue $ws.Range('E49') '  -1.60%  '
Set-TextValue $ws.Range('E50') '  -1.70%  '
Set-TextValue $ws.Range('D51') '0.0769'
Set-TextValue $ws.Range('E51') '  -1.01%  '
